# SendMoney_Schedule.xlsx update
#
# The "BEGIN UPDATE DC_FUND_TRANSFER_BENEFICIARY ..." cleanup SQL in B2
# hard-coded a literal account number ('06047900194203'). Replace that
# literal with the '{account_number}' template placeholder so the script
# works for whichever account the scenario actually used, instead of the
# one fixed account hard-coded before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldSql = $ws.Range("B2").Text
$newSql = $oldSql.Replace("'06047900194203'", "'{account_number}'")
$ws.Range("B2").Value = $newSql

# Move the sheet's selection/viewport from Q12 to B3, matching where the
# author left the cursor after making the edit above.
$ws.Range("B3").Select()
